{"js": "// Apply the translation-memory text corrections described by the diff.\n// Each entry: a search phrase (unique within the doc, or disambiguated by\n// style) and its replacement. We use Range.search + Range.insertText with\n// Word.InsertLocation.replace so the run keeps its original \"FuzzyMatch\"\n// character style.\nconst body = context.document.body;\n\nasync function replaceOnce(searchText, replacement, styleFilter) {\n  const results = body.search(searchText, { matchCase: true, matchWholeWord: false });\n  results.load(\"text,style\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    const item = results.items[i];\n    if (styleFilter && item.style !== styleFilter) {\n      continue;\n    }\n    item.insertText(replacement, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n// 1) \"... Mayo Clinic Remote Patient Monitoring Program\" -> \"... Mayo Clinic Lub Koom Kas Kev Saib Xyuas Tus Neeg Mob Nyob Deb\"\n//    (only the translated/\"FuzzyMatch\" cell, not the English source cell)\nawait replaceOnce(\n  \" Mayo Clinic Remote Patient Monitoring Program\",\n  \" Mayo Clinic Lub Koom Kas Kev Saib Xyuas Tus Neeg Mob Nyob Deb\",\n  \"FuzzyMatch\"\n);\n\n// 2) \"Qhov kev Pab Cuam Saib\" -> \"qhov Kev Saib\"\nawait replaceOnce(\"Qhov kev Pab Cuam Saib\", \"qhov Kev Saib\");\n\n// 3) \"Mob nyob deb ntawm lub\" -> \"Mob Nyob Deb ntawm lub\"\nawait replaceOnce(\"Mob nyob deb ntawm lub\", \"Mob Nyob Deb ntawm lub\");\n\n// 4) \"qab Koj Cov\" -> \"Qab Mus Rau Koj Cov\"\nawait replaceOnce(\"qab Koj Cov\", \"Qab Mus Rau Koj Cov\");\n\n// 5) \"yuav sai sai hu\" -> \"yuav hu\"\nawait replaceOnce(\"yuav sai sai hu\", \"yuav hu\");\n\n// 6) \"sijhawm\" -> \"sij hawm\" (word split, same visible text)\nawait replaceOnce(\"sijhawm\", \"sij hawm\");\n\n// 7) \"Pab Pawg Saib Xyuas Tus Neeg Mob Nyob deb ntawm\" -> \"peb pab pawg Kev Saib Xyuas Tus Neeg Mob Nyob Deb ntawm \"\nawait replaceOnce(\n  \"Pab Pawg Saib Xyuas Tus Neeg Mob Nyob deb ntawm\",\n  \"peb pab pawg Kev Saib Xyuas Tus Neeg Mob Nyob Deb ntawm \"\n);\n\n// 8) \"Pab Pawg Saib Xyuas Neeg Mob Nyob Hauv Koj Lub Chaw Saib Xyuas Neeg Mob\"\n//    -> \"Koj Pab Pawg Neeg Saib Xyuas Tus Neeg Mob Nyob Deb\"\nawait replaceOnce(\n  \"Pab Pawg Saib Xyuas Neeg Mob Nyob Hauv Koj Lub Chaw Saib Xyuas Neeg Mob\",\n  \"Koj Pab Pawg Neeg Saib Xyuas Tus Neeg Mob Nyob Deb\"\n);\n", "ps1": "# Apply the translation-memory text corrections described by the diff.\n# Uses Range.Find/Replacement (restricted to the \"FuzzyMatch\" character\n# style where needed) so only the translated/target table cells are\n# touched, never the English source cells that happen to share text.\n\nfunction Replace-FuzzyText {\n    param(\n        [__ComObject]$Doc,\n        [string]$FindText,\n        [string]$ReplaceText,\n        [string]$StyleName = $null\n    )\n\n    $range = $Doc.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Text = $FindText\n\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $ReplaceText\n\n    $useFormat = $false\n    if ($StyleName) {\n        $find.Style = $StyleName\n        $useFormat = $true\n    }\n\n    $find.Forward = $true\n    $find.Wrap = 0\n    $find.Format = $useFormat\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n\n    $find.Execute($FindText, $false, $false, $false, $false, $false, $true, 1, $useFormat, $ReplaceText, 2) | Out-Null\n}\n\n$d = $word.ActiveDocument\n\n# 1) \"... Mayo Clinic Remote Patient Monitoring Program\" -> \"... Mayo Clinic Lub Koom Kas Kev Saib Xyuas Tus Neeg Mob Nyob Deb\"\n#    (restrict to the \"FuzzyMatch\" style so the English source cell is left alone)\nReplace-FuzzyText $d \"Mayo Clinic Remote Patient Monitoring Program\" \"Mayo Clinic Lub Koom Kas Kev Saib Xyuas Tus Neeg Mob Nyob Deb\" \"FuzzyMatch\"\n\n# 2) \"Qhov kev Pab Cuam Saib\" -> \"qhov Kev Saib\"\nReplace-FuzzyText $d \"Qhov kev Pab Cuam Saib\" \"qhov Kev Saib\"\n\n# 3) \"Mob nyob deb ntawm lub\" -> \"Mob Nyob Deb ntawm lub\"\nReplace-FuzzyText $d \"Mob nyob deb ntawm lub\" \"Mob Nyob Deb ntawm lub\"\n\n# 4) \"qab Koj Cov\" -> \"Qab Mus Rau Koj Cov\"\nReplace-FuzzyText $d \"qab Koj Cov\" \"Qab Mus Rau Koj Cov\"\n\n# 5) \"yuav sai sai hu\" -> \"yuav hu\"\nReplace-FuzzyText $d \"yuav sai sai hu\" \"yuav hu\"\n\n# 6) \"sijhawm\" -> \"sij hawm\" (word split, same visible text)\nReplace-FuzzyText $d \"sijhawm\" \"sij hawm\"\n\n# 7) \"Pab Pawg Saib Xyuas Tus Neeg Mob Nyob deb ntawm\" -> \"peb pab pawg Kev Saib Xyuas Tus Neeg Mob Nyob Deb ntawm \"\nReplace-FuzzyText $d \"Pab Pawg Saib Xyuas Tus Neeg Mob Nyob deb ntawm\" \"peb pab pawg Kev Saib Xyuas Tus Neeg Mob Nyob Deb ntawm \"\n\n# 8) \"Pab Pawg Saib Xyuas Neeg Mob Nyob Hauv Koj Lub Chaw Saib Xyuas Neeg Mob\"\n#    -> \"Koj Pab Pawg Neeg Saib Xyuas Tus Neeg Mob Nyob Deb\"\nReplace-FuzzyText $d \"Pab Pawg Saib Xyuas Neeg Mob Nyob Hauv Koj Lub Chaw Saib Xyuas Neeg Mob\" \"Koj Pab Pawg Neeg Saib Xyuas Tus Neeg Mob Nyob Deb\"\n"}
